$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3844
$ws.Range("F5").Value = 1375
$ws.Range("F6").Value = 3840
$ws.Range("F8").Value = 205
$ws.Range("F9").Value = 58
$ws.Range("F10").Value = 8712
$ws.Range("F13").Value = 302
$ws.Range("F14").Value = 344
$ws.Range("F16").Value = 101
$ws.Range("F18").Value = 368
$ws.Range("F19").Value = 11038
$ws.Range("F25").Value = 42
$ws.Range("F28").Value = 192
$ws.Range("F30").Value = 161
$ws.Range("F34").Value = 2083
$ws.Range("F35").Value = 44
$ws.Range("F36").Value = 42
$ws.Range("F38").Value = 905
$ws.Range("F39").Value = 2570
$ws.Range("F40").Value = 286
$ws.Range("F41").Value = 2589
$ws.Range("F42").Value = 1250
$ws.Range("F43").Value = 176
$ws.Range("F44").Value = 760
$ws.Range("F46").Value = 349
$ws.Range("F47").Value = 49
$ws.Range("F48").Value = 91
$ws.Range("F49").Value = 85

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 7
$ws.Range("F9").Value = 9
$ws.Range("F21").Value = 31

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 35

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 3844
$ws.Range("F6").Value = 3840
$ws.Range("F10").Value = 205
$ws.Range("F11").Value = 58
$ws.Range("F12").Value = 8712
$ws.Range("F16").Value = 302
$ws.Range("F17").Value = 344
$ws.Range("F19").Value = 101
$ws.Range("F20").Value = 368
$ws.Range("F21").Value = 11038
$ws.Range("F23").Value = 42
$ws.Range("F26").Value = 192
$ws.Range("F28").Value = 161
$ws.Range("F30").Value = 2083
$ws.Range("F31").Value = 44
$ws.Range("F32").Value = 42
$ws.Range("F34").Value = 905
$ws.Range("F37").Value = 2570
$ws.Range("F38").Value = 286
$ws.Range("F39").Value = 2589
$ws.Range("F41").Value = 1250
$ws.Range("F42").Value = 176
$ws.Range("F43").Value = 760
$ws.Range("F45").Value = 349
$ws.Range("F47").Value = 49
$ws.Range("F48").Value = 91
$ws.Range("F49").Value = 85

